$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.402.75'
$ws.Range('E2').Value = '  -2.34%  '

$ws.Range('D3').Value = '3.698.43'
$ws.Range('E3').Value = '  -2.92%  '

$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').Value = '''691.83'
$ws.Range('E5').Value = '  -0.92%  '

$ws.Range('D6').Value = '''163.04'
$ws.Range('E6').Value = '  -5.12%  '

$ws.Range('D7').Value = '3.697.40'
$ws.Range('E7').Value = '  -2.94%  '

$ws.Range('E8').Value = '  +0.09%  '

$ws.Range('E9').Value = '  -4.70%  '

$ws.Range('E10').Value = '  -8.14%  '

$ws.Range('D11').Value = '''7.42'
$ws.Range('E11').Value = '  -1.32%  '

$ws.Range('E12').Value = '  -4.32%  '

$ws.Range('D13').Value = '''0.0000239'
$ws.Range('E13').Value = '  -5.24%  '

$ws.Range('D14').Value = '''33.43'
$ws.Range('E14').Value = '  -7.11%  '

$ws.Range('D15').Value = '4.322.02'
$ws.Range('E15').Value = '  -3.04%  '

$ws.Range('D16').Value = '3.704.66'
$ws.Range('E16').Value = '  -2.82%  '

$ws.Range('D17').Value = '69.472.89'
$ws.Range('E17').Value = '  -2.24%  '

$ws.Range('E18').Value = '  -0.86%  '

$ws.Range('D19').Value = '''16.26'
$ws.Range('E19').Value = '  -6.93%  '

$ws.Range('D20').Value = '''6.60'
$ws.Range('E20').Value = '  -7.58%  '

$ws.Range('D21').Value = '''481.05'
$ws.Range('E21').Value = '  -6.29%  '

$ws.Range('D22').Value = '''9.98'
$ws.Range('E22').Value = '  -6.53%  '

$ws.Range('D23').Value = '''0.664'
$ws.Range('E23').Value = '  -7.15%  '

$ws.Range('D25').Value = '3.846.37'
$ws.Range('E25').Value = '  -2.95%  '

$ws.Range('E26').Value = '  -9.33%  '

$ws.Range('E27').Value = '  -0.01%  '

$ws.Range('D28').Value = '''11.40'
$ws.Range('E28').Value = '  -5.40%  '

$ws.Range('D29').Value = '''9.52'
$ws.Range('E29').Value = '  -8.41%  '

$ws.Range('E30').Value = '  -10.41%  '

$ws.Range('D31').Value = '''2.72'
$ws.Range('E31').Value = '  -9.90%  '

$ws.Range('D32').Value = '''6.85'
$ws.Range('E32').Value = '  -7.71%  '

$ws.Range('D33').Value = '''2.07'
$ws.Range('E33').Value = '  -7.73%  '

$ws.Range('B34').Value = 'Kaspa'
$ws.Range('C34').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D34').Value = '''0.168'
$ws.Range('E34').Value = '  -4.75%  '

$ws.Range('B35').Value = 'EthereumClassic'
$ws.Range('C35').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D35').Value = '''27.04'
$ws.Range('E35').Value = '  -7.01%  '

$ws.Range('B36').Value = 'Binance-PegBSC-USD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D36').Value = '''1.00'
$ws.Range('E36').Value = '  +0.20%  '

$ws.Range('D37').Value = '3.667.82'
$ws.Range('E37').Value = '  -2.77%  '

$ws.Range('D38').Value = '''8.49'
$ws.Range('E38').Value = '  -7.25%  '

$ws.Range('D39').Value = '''6.38'
$ws.Range('E39').Value = '  +6.24%  '

$ws.Range('E40').Value = '  -2.89%  '

$ws.Range('D41').Value = '''0.0932'
$ws.Range('E41').Value = '  -7.77%  '

$ws.Range('E43').Value = '  +0.00%  '

$ws.Range('E44').Value = '  -6.43%  '

$ws.Range('D45').Value = '''163.76'
$ws.Range('E45').Value = '  -5.58%  '

$ws.Range('D46').Value = '''48.10'
$ws.Range('E46').Value = '  -2.67%  '

$ws.Range('D47').Value = '''30.14'
$ws.Range('E47').Value = '  +1.59%  '

$ws.Range('D48').Value = '''2.81'
$ws.Range('E48').Value = '  -15.25%  '

$ws.Range('B49').Value = 'SuiNetwork'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D49').Value = '''1.15'
$ws.Range('E49').Value = '  -1.18%  '

$ws.Range('B50').Value = 'ONDO'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D50').Value = '''1.34'
$ws.Range('E50').Value = '  -1.51%  '

$ws.Range('D51').Value = '''0.000285'
$ws.Range('E51').Value = '  -8.15%  '
